$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "IDE + CTR (ctr$ hit)" row (old row 4) entirely; the row
# below it ("IDE + CTR (ctr$ miss)") shifts up to become row 4.
$ws.Rows.Item(4).Delete()

# Remove the trailing "ther" column (old column G) entirely.
$ws.Columns.Item(7).Delete()

# Rename the header labels.
$ws.Range("B1").Value = "GCM"
$ws.Range("C1").Value = "XTS"
$ws.Range("D1").Value = "otp"
$ws.Range("E1").Value = "counteraccess"
$ws.Range("F1").Value = "cxl"

# Row 2 - Adaptive
$ws.Range("B2").Value = 0.01090218641
$ws.Range("C2").Value = 0.09427379885
$ws.Range("D2").Value = 0.09811967769
$ws.Range("E2").Value = 0.04360874564
$ws.Range("F2").Value = 0.01817031068

# Row 3 - IDE + XTS
$ws.Range("B3").Value = 0.654764097
$ws.Range("C3").Value = 0.3270655923
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.01817031068

# Row 4 - IDE + CTR (ctr$ miss)  [was row 5 before the row delete above]
$ws.Range("B4").Value = 0.654764097
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0.1505949739
$ws.Range("E4").Value = 0.4000215328
$ws.Range("F4").Value = 0.01817031068
